$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Update the value that changed in row 3
$ws.Range("B3").Value = 285317708

# Extend the demand profile with a full 12-timestep (+header) series: rows 4-14,
# A = timestep index (2..12), B = constant demand value.
$row = 4
for ($t = 2; $t -le 12; $t++) {
    $ws.Cells.Item($row, 1).Value = $t
    $ws.Cells.Item($row, 2).Value = 285317708
    $row++
}

# Column B widened to fit the longer values
$ws.Columns.Item(2).ColumnWidth = 10.14

# Make Demand the active sheet/tab and leave the selection where the author left it
$ws.Select()
$ws.Range("C13").Select()
